$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Row 42
$ws.Range("A42").Value = "VerifyPostRecordDetails"
$ws.Range("B42").Value = "OPQA-370"
$ws.Range("C42").Value = "Verify that user contributed articles display the information about the author"
$ws.Range("D42").Value = "Y"
$ws.Range("E42").Value = "PASS"

# Row 43
$ws.Range("A43").Value = "SeacrhAndViewOwnPost"
$ws.Range("B43").Value = "OPQA-415"
$ws.Range("C43").Value = "Verify that user is able to search the  posts a user authored themselves and view them."
$ws.Range("D43").Value = "Y"
$ws.Range("E43").Value = "PASS"

# Row 44
$ws.Range("A44").Value = "SeacrhAndViewOthersPost"
$ws.Range("B44").Value = "OPQA-416"
$ws.Range("C44").Value = "Verify that user is able to search the posts of others and view them."
$ws.Range("D44").Value = "Y"
$ws.Range("E44").Value = "PASS"

# Copy formatting for the A/C/D/E columns from an existing row that uses
# the plain bordered style (s=1) on another sheet, and the hyperlink-like
# style (s=19) used for column B from this sheet's own B2 cell -- this
# matches the styling used by the appended rows in the target workbook.
$wsSrc = $wb.Worksheets.Item("CommentsProfanityWordsCheckTest")

$wsSrc.Range("A2").Copy()
$ws.Range("A42:A44").PasteSpecial(-4122)
$ws.Range("C42:C44").PasteSpecial(-4122)
$ws.Range("D42:D44").PasteSpecial(-4122)
$ws.Range("E42:E44").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("B42:B44").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Update the view: scroll so row 28 is at the top and select B44 (the last
# edited cell), matching the saved sheet view state.
$ws.Application.Goto($ws.Range("A28"))
$ws.Range("B44").Select()
